$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25; everything from row 25 downward shifts down by one.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with a fresh weekly entry (same series as the
# last existing record, but for a new date).
$ws.Cells.Item(25, 1).Value = 7
$ws.Cells.Item(25, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(25, 3).Value = "Ñuble"
$ws.Cells.Item(25, 4).Value = 45030
$ws.Cells.Item(25, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25, 5).Value = 16
$ws.Cells.Item(25, 6).Value = 100112030
$ws.Cells.Item(25, 7).Value = "Poroto granado"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 60
$ws.Cells.Item(25, 11).Value = 35000
$ws.Cells.Item(25, 12).Value = 35000
$ws.Cells.Item(25, 13).Value = 35000
$ws.Cells.Item(25, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 1400
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
